# "Add Tools 1 and 2"
#
# Concentration_Time_Data: populate the "Tool 2" column (M) with the
# refreshed TCE readings (mirrors the "Tool 1" column L / raw reading
# column E for each event row).
#
# Monitoring_Well_Information: update the monitoring-well GPS
# coordinates (columns B/C) with newly surveyed lat/long readings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Concentration_Time_Data  (column M = "Tool 2")
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Concentration_Time_Data")
[void]$ws2.Activate()

# row -> new M-column value ($null clears the cell)
$mValues = [ordered]@{
    3  = 37.1
    4  = 41.9
    5  = 13
    6  = 5.0999999999999996
    7  = 11.5
    8  = 5
    9  = $null
    10 = 4.5999999999999996
    11 = 1.85
    12 = $null
    13 = 1.8
    14 = 1.2
    15 = 1
    16 = 1.2
    17 = 1
    18 = 0.8
    20 = 0.7
    21 = 0.5
}

foreach ($row in $mValues.Keys) {
    $val = $mValues[$row]
    $cell = $ws2.Range("M$row")
    if ($null -eq $val) {
        [void]$cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}

[void]$ws2.Range("I3:I22").Select()

# ---------------------------------------------------------------------
# Sheet: Monitoring_Well_Information  (refreshed lat/long readings)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Monitoring_Well_Information")
[void]$ws3.Activate()

# row -> @(latitude, longitude)
$coords = [ordered]@{
    2  = @(29.731660000000002, -95.412599999999998)
    3  = @(29.732970000000002, -95.413979999999995)
    4  = @(29.733080000000001, -95.413020000000003)
    5  = @(29.732679999999998, -95.411810000000003)
    6  = @(29.733730000000001, -95.412239999999997)
    7  = @(29.733250000000002, -95.411389999999997)
    8  = @(29.733989999999999, -95.410929999999993)
    9  = @(29.732690000000002, -95.413179999999997)
    10 = @(29.732559999999999, -95.412350000000004)
}

foreach ($row in $coords.Keys) {
    $pair = $coords[$row]
    $ws3.Range("B$row").Value = $pair[0]
    $ws3.Range("C$row").Value = $pair[1]
}

[void]$ws3.Range("E15").Select()
